$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.225.12"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.295.52"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'533.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.41%  "
$ws.Range("D6").Value = "'131.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +2.61%  "
$ws.Range("D9").Value = "2.293.71"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("E10").Value = "  -2.99%  "
$ws.Range("D11").Value = "'5.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").Value = "'23.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "2.706.68"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("D16").Value = "58.149.28"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "2.321.73"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "'4.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.63%  "
$ws.Range("D21").Value = "'312.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'8.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("D29").Value = "'169.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("D31").Value = "0.0₃0721"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("E34").Value = "  -3.76%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("D40").Value = "'38.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("D42").Value = "'141.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").Value = "'288.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.08%  "
$ws.Range("D44").Value = "'3.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "'0.0494"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "'0.556"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "'18.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.24%  "
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").Value = "'10.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("E51").Value = "  -0.68%  "
